$d = $word.ActiveDocument

# --- Step 1: modifications that do not change the paragraph count ---

# Paragraph 9 (numId=2): "Создание стилей CSS" -> remove the "Создание стилей "
# run and change the remaining run's text from "CSS" to "Test-cases",
# keeping that run's own formatting (en-US lang) intact.
$p9 = $d.Paragraphs.Item(9)
$p9Start = $p9.Range.Start
$toRemove = $d.Range($p9Start, $p9Start + 16)
if ($toRemove.Text -eq "Создание стилей ") {
    $toRemove.Delete()
}
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Find.Execute("CSS", $true, $false, $false, $false, $false, $true, 1, $false, "Test-cases", 2) | Out-Null

# Paragraph 13 (numId=3, Roman Nikulin section): "Диаграмма классов" -> "Диаграмма последовательности"
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Find.Execute("Диаграмма классов", $true, $false, $false, $false, $false, $true, 1, $false, "Диаграмма последовательности", 2) | Out-Null

# --- Step 2: whole-paragraph deletions (highest index first so indices stay valid) ---

# 18: "Создание стартовой страницы сайта" (numId=3)
$d.Paragraphs.Item(18).Range.Delete()

# 17: "Создание модуля авторизации" (numId=3)
$d.Paragraphs.Item(17).Range.Delete()

# 16: "IDEF0 диаграммы для модулей системы" (numId=3)
$d.Paragraphs.Item(16).Range.Delete()

# 14: "Диаграмма последовательности" (duplicate, numId=3)
$d.Paragraphs.Item(14).Range.Delete()

# 11: "IDEF" + "0 " + "диаграммы для модулей системы" (numId=2)
$d.Paragraphs.Item(11).Range.Delete()

# 10: "Test-cases" (duplicate, numId=2)
$d.Paragraphs.Item(10).Range.Delete()

# 6: "IDEF0 диаграммы для модулей системы" (numId=1)
$d.Paragraphs.Item(6).Range.Delete()

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"
